{"js": "// Auto-generated replacements: update math problem cells in the table\n// to match the new arithmetic exercises from the diff.\nconst replacements = [\n  [\"59+23=82\", \"96-76=20\"],\n  [\"54+45=99\", \"27+13=40\"],\n  [\"9+74=83\", \"64-45=19\"],\n  [\"75-63=12\", \"8-1=7\"],\n  [\"5+47=52\", \"3+0=3\"],\n  [\"45-37=8\", \"49-16=33\"],\n  [\"67-52=15\", \"60-39=21\"],\n  [\"18-15=3\", \"45+41=86\"],\n  [\"32-18=14\", \"84+7=91\"],\n  [\"17+72=89\", \"18+7=25\"],\n  [\"33-11=22\", \"35-3=32\"],\n  [\"94-71=23\", \"67-17=50\"],\n  [\"51-36=15\", \"62-48=14\"],\n  [\"66+8=74\", \"19+47=66\"],\n  [\"19+28=47\", \"81-20=61\"],\n  [\"6+16=22\", \"26-9=17\"],\n  [\"75-61=14\", \"90-1=89\"],\n  [\"72-7=65\", \"69+9=78\"],\n  [\"7+1=8\", \"64+17=81\"],\n  [\"51+24=75\", \"91-64=27\"],\n  [\"85-13=72\", \"24-12=12\"],\n  [\"36+38=74\", \"4+90=94\"],\n  [\"98-71=27\", \"58-24=34\"],\n  [\"34-29=5\", \"66-38=28\"],\n  [\"25+14=39\", \"14+38=52\"],\n  [\"34-28=6\", \"13+34=47\"],\n  [\"33+6=39\", \"46+53=99\"],\n  [\"27+4=31\", \"16+54=70\"],\n  [\"39+49=88\", \"52+45=97\"],\n  [\"19+23=42\", \"94-42=52\"],\n  [\"9+33=42\", \"54-23=31\"],\n  [\"31+52=83\", \"87-58=29\"],\n  [\"38+6=44\", \"64+2=66\"],\n  [\"72-28=44\", \"67-34=33\"],\n  [\"63-36=27\", \"21+17=38\"],\n  [\"76-11=65\", \"43-7=36\"],\n  [\"0+48=48\", \"61-54=7\"],\n  [\"58+21=79\", \"94-76=18\"],\n  [\"49-5=44\", \"58+7=65\"],\n  [\"81+16=97\", \"52-25=27\"],\n  [\"82-51=31\", \"41-36=5\"],\n  [\"61-52=9\", \"79-44=35\"],\n  [\"56-46=10\", \"11-5=6\"],\n  [\"19+29=48\", \"75-16=59\"],\n  [\"61-43=18\", \"42+14=56\"],\n  [\"78-46=32\", \"23+74=97\"],\n  [\"69-16=53\", \"28+9=37\"],\n  [\"47+17=64\", \"64-37=27\"],\n  [\"64+9=73\", \"58-40=18\"],\n  [\"1+57=58\", \"37-33=4\"],\n  [\"73-64=9\", \"28-16=12\"],\n  [\"38+51=89\", \"66-6=60\"],\n  [\"0+79=79\", \"46+3=49\"],\n  [\"84-39=45\", \"35-2=33\"],\n  [\"0+63=63\", \"80+9=89\"],\n  [\"4+95=99\", \"45+44=89\"],\n  [\"1+27=28\", \"56+42=98\"],\n  [\"77+0=77\", \"84-67=17\"],\n  [\"21+14=35\", \"72+19=91\"],\n  [\"31+22=53\", \"57+0=57\"],\n  [\"87-81=6\", \"40+45=85\"],\n  [\"20+73=93\", \"96-89=7\"],\n  [\"52-35=17\", \"65-39=26\"],\n  [\"71-70=1\", \"24-3=21\"],\n  [\"0+65=65\", \"40-19=21\"],\n  [\"57+19=76\", \"59+18=77\"],\n  [\"6+29=35\", \"91-31=60\"],\n  [\"61-2=59\", \"97-31=66\"],\n  [\"42+2=44\", \"30-26=4\"],\n  [\"21+16=37\", \"49+36=85\"],\n  [\"71-47=24\", \"45-16=29\"],\n  [\"46+27=73\", \"59-1=58\"],\n  [\"77+14=91\", \"59-11=48\"],\n  [\"95-87=8\", \"47-34=13\"],\n  [\"13+24=37\", \"22+67=89\"],\n  [\"70-62=8\", \"75-35=40\"],\n  [\"35+39=74\", \"49-11=38\"],\n  [\"85-26=59\", \"17+60=77\"],\n  [\"92-29=63\", \"84+13=97\"],\n  [\"68-58=10\", \"14+74=88\"],\n  [\"89-6=83\", \"44-31=13\"],\n  [\"5+51=56\", \"6+15=21\"],\n  [\"64-42=22\", \"15+63=78\"],\n  [\"11-4=7\", \"95-61=34\"],\n  [\"32-21=11\", \"40+14=54\"],\n  [\"3+9=12\", \"19+28=47\"],\n  [\"66-14=52\", \"39+10=49\"],\n  [\"35-20=15\", \"47-13=34\"],\n  [\"76+5=81\", \"41+19=60\"],\n  [\"54+36=90\", \"8+39=47\"],\n  [\"40-11=29\", \"43-17=26\"],\n  [\"38+10=48\", \"30+23=53\"],\n  [\"6+52=58\", \"30-11=19\"],\n  [\"39+4=43\", \"63-59=4\"],\n  [\"35-27=8\", \"45+49=94\"],\n  [\"59-18=41\", \"3+60=63\"],\n  [\"82-68=14\", \"49-36=13\"],\n  [\"3+65=68\", \"46+11=57\"],\n  [\"24-24=0\", \"82+5=87\"],\n  [\"79+15=94\", \"73-2=71\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Auto-generated replacements: update math problem cells in the table\n# to match the new arithmetic exercises from the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"59+23=82\", \"96-76=20\"),\n    @(\"54+45=99\", \"27+13=40\"),\n    @(\"9+74=83\", \"64-45=19\"),\n    @(\"75-63=12\", \"8-1=7\"),\n    @(\"5+47=52\", \"3+0=3\"),\n    @(\"45-37=8\", \"49-16=33\"),\n    @(\"67-52=15\", \"60-39=21\"),\n    @(\"18-15=3\", \"45+41=86\"),\n    @(\"32-18=14\", \"84+7=91\"),\n    @(\"17+72=89\", \"18+7=25\"),\n    @(\"33-11=22\", \"35-3=32\"),\n    @(\"94-71=23\", \"67-17=50\"),\n    @(\"51-36=15\", \"62-48=14\"),\n    @(\"66+8=74\", \"19+47=66\"),\n    @(\"19+28=47\", \"81-20=61\"),\n    @(\"6+16=22\", \"26-9=17\"),\n    @(\"75-61=14\", \"90-1=89\"),\n    @(\"72-7=65\", \"69+9=78\"),\n    @(\"7+1=8\", \"64+17=81\"),\n    @(\"51+24=75\", \"91-64=27\"),\n    @(\"85-13=72\", \"24-12=12\"),\n    @(\"36+38=74\", \"4+90=94\"),\n    @(\"98-71=27\", \"58-24=34\"),\n    @(\"34-29=5\", \"66-38=28\"),\n    @(\"25+14=39\", \"14+38=52\"),\n    @(\"34-28=6\", \"13+34=47\"),\n    @(\"33+6=39\", \"46+53=99\"),\n    @(\"27+4=31\", \"16+54=70\"),\n    @(\"39+49=88\", \"52+45=97\"),\n    @(\"19+23=42\", \"94-42=52\"),\n    @(\"9+33=42\", \"54-23=31\"),\n    @(\"31+52=83\", \"87-58=29\"),\n    @(\"38+6=44\", \"64+2=66\"),\n    @(\"72-28=44\", \"67-34=33\"),\n    @(\"63-36=27\", \"21+17=38\"),\n    @(\"76-11=65\", \"43-7=36\"),\n    @(\"0+48=48\", \"61-54=7\"),\n    @(\"58+21=79\", \"94-76=18\"),\n    @(\"49-5=44\", \"58+7=65\"),\n    @(\"81+16=97\", \"52-25=27\"),\n    @(\"82-51=31\", \"41-36=5\"),\n    @(\"61-52=9\", \"79-44=35\"),\n    @(\"56-46=10\", \"11-5=6\"),\n    @(\"19+29=48\", \"75-16=59\"),\n    @(\"61-43=18\", \"42+14=56\"),\n    @(\"78-46=32\", \"23+74=97\"),\n    @(\"69-16=53\", \"28+9=37\"),\n    @(\"47+17=64\", \"64-37=27\"),\n    @(\"64+9=73\", \"58-40=18\"),\n    @(\"1+57=58\", \"37-33=4\"),\n    @(\"73-64=9\", \"28-16=12\"),\n    @(\"38+51=89\", \"66-6=60\"),\n    @(\"0+79=79\", \"46+3=49\"),\n    @(\"84-39=45\", \"35-2=33\"),\n    @(\"0+63=63\", \"80+9=89\"),\n    @(\"4+95=99\", \"45+44=89\"),\n    @(\"1+27=28\", \"56+42=98\"),\n    @(\"77+0=77\", \"84-67=17\"),\n    @(\"21+14=35\", \"72+19=91\"),\n    @(\"31+22=53\", \"57+0=57\"),\n    @(\"87-81=6\", \"40+45=85\"),\n    @(\"20+73=93\", \"96-89=7\"),\n    @(\"52-35=17\", \"65-39=26\"),\n    @(\"71-70=1\", \"24-3=21\"),\n    @(\"0+65=65\", \"40-19=21\"),\n    @(\"57+19=76\", \"59+18=77\"),\n    @(\"6+29=35\", \"91-31=60\"),\n    @(\"61-2=59\", \"97-31=66\"),\n    @(\"42+2=44\", \"30-26=4\"),\n    @(\"21+16=37\", \"49+36=85\"),\n    @(\"71-47=24\", \"45-16=29\"),\n    @(\"46+27=73\", \"59-1=58\"),\n    @(\"77+14=91\", \"59-11=48\"),\n    @(\"95-87=8\", \"47-34=13\"),\n    @(\"13+24=37\", \"22+67=89\"),\n    @(\"70-62=8\", \"75-35=40\"),\n    @(\"35+39=74\", \"49-11=38\"),\n    @(\"85-26=59\", \"17+60=77\"),\n    @(\"92-29=63\", \"84+13=97\"),\n    @(\"68-58=10\", \"14+74=88\"),\n    @(\"89-6=83\", \"44-31=13\"),\n    @(\"5+51=56\", \"6+15=21\"),\n    @(\"64-42=22\", \"15+63=78\"),\n    @(\"11-4=7\", \"95-61=34\"),\n    @(\"32-21=11\", \"40+14=54\"),\n    @(\"3+9=12\", \"19+28=47\"),\n    @(\"66-14=52\", \"39+10=49\"),\n    @(\"35-20=15\", \"47-13=34\"),\n    @(\"76+5=81\", \"41+19=60\"),\n    @(\"54+36=90\", \"8+39=47\"),\n    @(\"40-11=29\", \"43-17=26\"),\n    @(\"38+10=48\", \"30+23=53\"),\n    @(\"6+52=58\", \"30-11=19\"),\n    @(\"39+4=43\", \"63-59=4\"),\n    @(\"35-27=8\", \"45+49=94\"),\n    @(\"59-18=41\", \"3+60=63\"),\n    @(\"82-68=14\", \"49-36=13\"),\n    @(\"3+65=68\", \"46+11=57\"),\n    @(\"24-24=0\", \"82+5=87\"),\n    @(\"79+15=94\", \"73-2=71\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 0)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n    $rng.Text = $newText\n}\n"}
